# "Season up to 1/17"
# The game that was next up (2024-01-15 @ ATL, the first row of the "Next"
# schedule sheet) has now been played. Record its final box-score line as
# the new last row of the "Games" results sheet, and drop it from the
# "Next" schedule (which shifts every remaining game up by one row, with
# no new game appended at the bottom).

$wb = $excel.ActiveWorkbook
$gamesWs = $wb.Worksheets.Item("Games")
$nextWs  = $wb.Worksheets.Item("Next")

# --- 1. Append the newly-played game to the "Games" sheet ---------------
$newRow = $gamesWs.UsedRange.Rows.Count + 1

$gamesWs.Cells.Item($newRow, 1).Value  = 39        # Game
$gamesWs.Cells.Item($newRow, 2).Value  = 45306     # Date
$gamesWs.Cells.Item($newRow, 2).NumberFormat = $gamesWs.Cells.Item($newRow - 1, 2).NumberFormat
$gamesWs.Cells.Item($newRow, 3).Value  = -2        # Streak
$gamesWs.Cells.Item($newRow, 4).Value  = 99        # Pts
$gamesWs.Cells.Item($newRow, 5).Value  = 105.1     # Pace
$gamesWs.Cells.Item($newRow, 6).Value  = 0.458     # eFG
$gamesWs.Cells.Item($newRow, 7).Value  = 14.2      # TOV
$gamesWs.Cells.Item($newRow, 8).Value  = 17.3      # ORB
$gamesWs.Cells.Item($newRow, 9).Value  = 0.126     # FTR
$gamesWs.Cells.Item($newRow, 10).Value = 94.2      # ORT
$gamesWs.Cells.Item($newRow, 11).Value = "ATL"     # OppID
$gamesWs.Cells.Item($newRow, 12).Value = 109       # OppPts
$gamesWs.Cells.Item($newRow, 13).Value = 0.5       # OppeFG
$gamesWs.Cells.Item($newRow, 14).Value = 13.5      # OppTOV
$gamesWs.Cells.Item($newRow, 15).Value = 27.1      # OppORB
$gamesWs.Cells.Item($newRow, 16).Value = 0.172     # OppFTR
$gamesWs.Cells.Item($newRow, 17).Value = 103.7     # OppORT
$gamesWs.Cells.Item($newRow, 18).Value = 0         # Location
$gamesWs.Cells.Item($newRow, 19).Value = 0         # Target

# --- 2. Remove that game from the "Next" schedule ------------------------
# Row 2 is the ATL game on 45306 that has now been played; deleting it
# shifts every following row up by one (and drops the old last row).
$nextWs.Rows.Item(2).Delete()
